$d = $word.ActiveDocument

function Get-BookmarkByName($doc, $name) {
    $bms = $doc.Bookmarks
    for ($i = 1; $i -le $bms.Count; $i++) {
        $b = $bms.Item($i)
        if ($b.Name -eq $name) { return $b }
    }
    return $null
}

# 1. Update the report date from 03 to 05 September 2021 (the "Date" styled paragraph).
$dateParaRange = $d.Paragraphs(3).Range
$dateParaRange.Find.Execute("03 September, 2021", $true, $false, $false, $false, $false, $true, 1, $false, "05 September, 2021", 2) | Out-Null

# 2. Re-create the "tab:OverviewTable" bookmark on the same range so it gets a fresh id.
$bm1 = Get-BookmarkByName $d "tab:OverviewTable"
$bm1Start = $bm1.Start
$bm1End = $bm1.End
$bm1Range = $d.Range($bm1Start, $bm1End)
$d.Bookmarks.Add("tab:OverviewTable", $bm1Range)

# 3. Re-create the "tab:StateLevelTable" bookmark on the same range so it gets a fresh id.
$bm2 = Get-BookmarkByName $d "tab:StateLevelTable"
$bm2Start = $bm2.Start
$bm2End = $bm2.End
$bm2Range = $d.Range($bm2Start, $bm2End)
$d.Bookmarks.Add("tab:StateLevelTable", $bm2Range)

Write-Host "done"
